$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the selection (active cell) on the sheet view
$ws.Range("H3").Select()

# Update H2 formula: now adds 64 to B2
$ws.Range("H2").Formula = '=$B$2+64'
